$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
Write-Host $shape.Name
Write-Host $shape.HasTable
$tbl = $shape.Table
Write-Host $tbl.Style
$tbl.ApplyStyle("{DF9E3A07-404A-4632-8D5F-6B45F7384861}")
Write-Host $tbl.Style
